$wb = $excel.ActiveWorkbook

# --- Hoja1: update the conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 5.36 = 21447.72 pesos`n✅ 21447.72 pesos = 5.35 = 948.71 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 186.5
$ws2.Range("O10").Value = 4000
$ws2.Range("N12").Value = 4008.5
$ws2.Range("O12").Value = 177.31
